$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.807.08"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.099.28"
$ws.Range("E3").Value = "  +3.78%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "387.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0864"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "3.589.73"
$ws.Range("E13").Value = "  +3.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "3.105.14"
$ws.Range("E16").Value = "  +3.92%  "
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.73%  "
$ws.Range("D19").Value = "51.932.65"
$ws.Range("E20").Value = "  +2.49%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.41%  "
$ws.Range("E27").Value = "  +3.83%  "
$ws.Range("E28").Value = "  +2.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.45%  "
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0451"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.296"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.66%  "
$ws.Range("E40").Value = "  +2.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.63%  "
$ws.Range("E47").Value = "  +4.12%  "
$ws.Range("E48").Value = "  +2.52%  "
$ws.Range("D49").Value = "2.050.44"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").Value = "3.410.08"
$ws.Range("E50").Value = "  +3.93%  "
$ws.Range("E51").Value = "  +7.18%  "
